$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "code"
$ws.Range("C1").Value = "exchangeRate"
$ws.Range("D1").Value = "symbol"
$ws.Range("E1").Value = "paymentMethods"

$ws.Range("K6").Select()
